$p = $ppt.ActivePresentation
Write-Host "HasNotesMaster: $($p.HasNotesMaster)"
Write-Host "HasHandoutMaster: $($p.HasHandoutMaster)"
$hm = $p.HandoutMaster
Write-Host "Handout: $hm"
try {
  $cs = $hm.ColorScheme
  for ($i=1;$i -le $cs.Count;$i++){
    Write-Host "$i : $($cs.Colors($i).RGB)"
  }
} catch { Write-Host "ERR: $_" }
